$d = $word.ActiveDocument

# The paragraph that contains the "10mathew10@gmail.com" hyperlink is the
# last of the four "Team members' emails" entries (currently paragraph 5).
# Everything after it (Support email section, Libraries used section,
# Resources section, Articles section) is being removed, except that the
# "_GoBack" bookmark that used to sit at the very end of the document is
# kept, and now ends up anchored right after that hyperlink.

# Locate the "10mathew10@gmail.com" hyperlink so we can compute a reliable
# (non hard-coded) insertion point for the bookmark once the trailing
# content has been removed.
$lastHyperlink = $d.Hyperlinks.Item(4)
$bmPos = $lastHyperlink.Range.End - 1
$bmRange = $d.Range($bmPos, $bmPos)

# Paragraph 5 is the one ending with the "10mathew10@gmail.com" hyperlink;
# delete everything from the start of paragraph 6 through the end of the
# document (this removes paragraphs 6-24 in one aligned operation).
$deleteStart = $d.Paragraphs.Item(6).Range.Start
$deleteRange = $d.Range($deleteStart, $d.Content.End)
$deleteRange.Delete()

# Re-create the "_GoBack" bookmark right after the hyperlink, now the very
# end of the document content.
$d.Bookmarks.Add("_GoBack", $bmRange)
